$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 92
$ws.Range("I11").Value = 92
$ws.Range("K11").Value = 92
$ws.Range("M11").Value = 48
$ws.Range("H15").Value = 2091.745
$ws.Range("I15").Value = 2091.745
$ws.Range("K15").Value = 6275.235
$ws.Range("M15").Value = -6106.235
$ws.Range("H39").Value = 3539
$ws.Range("I39").Value = 80
$ws.Range("K39").Value = 240
$ws.Range("M39").Value = 56
$ws.Range("H41").Value = 2769.2856
$ws.Range("J41").Value = 3476.8
$ws.Range("L41").Value = 3476.8
$ws.Range("N41").Value = -4356.8
$ws.Range("H48").Value = 6500
$ws.Range("J48").Value = 6500
$ws.Range("L48").Value = 19500
$ws.Range("N48").Value = -20084
$ws.Range("H56").Value = 6500
$ws.Range("J56").Value = 6500
$ws.Range("L56").Value = 19500
$ws.Range("N56").Value = -20568
$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36248
$ws.Range("H64").Value = 29109086
$ws.Range("I64").Value = 27548486
$ws.Range("J64").Value = 31254908
$ws.Range("K64").Value = 27548486
$ws.Range("L64").Value = 31254908
$ws.Range("M64").Value = -27548238
$ws.Range("N64").Value = -31255404
$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 35000
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111240
$ws.Range("H67").Value = 29109086
$ws.Range("I67").Value = 27548486
$ws.Range("J67").Value = 31254908
$ws.Range("K67").Value = 27548486
$ws.Range("L67").Value = 31254908
$ws.Range("M67").Value = -27547628
$ws.Range("N67").Value = -31256624
$ws.Range("H68").Value = 66990
$ws.Range("J68").Value = 66990
$ws.Range("L68").Value = 66990
$ws.Range("N68").Value = -68488
$ws.Range("H69").Value = 5901
$ws.Range("I69").Value = 5901
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 17703
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -16829
$ws.Range("N69").ClearContents()
$ws.Range("H71").Value = 66990
$ws.Range("J71").Value = 66990
$ws.Range("L71").Value = 200970
$ws.Range("N71").Value = -208458
$ws.Range("H72").Value = 5901
$ws.Range("I72").Value = 5901
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 53109
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -48741
$ws.Range("N72").ClearContents()
$ws.Range("H76").Value = 3974.5
$ws.Range("I76").Value = 3459.6
$ws.Range("K76").Value = 3459.6
$ws.Range("M76").Value = -3144.6
$ws.Range("H79").Value = 3974.5
$ws.Range("I79").Value = 3459.6
$ws.Range("K79").Value = 3459.6
$ws.Range("M79").Value = -2367.6
$ws.Range("H80").Value = 25000390
$ws.Range("I80").Value = 50000276
$ws.Range("J80").Value = 504
$ws.Range("K80").Value = 150000828
$ws.Range("L80").Value = 1512
$ws.Range("M80").Value = -149999830
$ws.Range("N80").Value = -3508
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H82").Value = 7521.273
$ws.Range("I82").Value = 4586.8
$ws.Range("K82").Value = 13760.4
$ws.Range("M82").Value = -13354.4
$ws.Range("H83").Value = 25000390
$ws.Range("I83").Value = 50000276
$ws.Range("J83").Value = 504
$ws.Range("K83").Value = 450002484
$ws.Range("L83").Value = 4536
$ws.Range("M83").Value = -449997492
$ws.Range("N83").Value = -14520
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H85").Value = 7521.273
$ws.Range("I85").Value = 4586.8
$ws.Range("K85").Value = 13760.4
$ws.Range("M85").Value = -12356.4
$ws.Range("H86").Value = 64517490
$ws.Range("I86").Value = 86022960
$ws.Range("J86").Value = 1075
$ws.Range("K86").Value = 86022960
$ws.Range("L86").Value = 1075
$ws.Range("M86").Value = -86021837
$ws.Range("N86").Value = -3321
$ws.Range("H87").Value = 82700
$ws.Range("J87").Value = 82700
$ws.Range("L87").Value = 82700
$ws.Range("N87").Value = -85196
$ws.Range("H89").Value = 64517490
$ws.Range("I89").Value = 86022960
$ws.Range("J89").Value = 1075
$ws.Range("K89").Value = 430114800
$ws.Range("L89").Value = 5375
$ws.Range("M89").Value = -430109184
$ws.Range("N89").Value = -16607
$ws.Range("H90").Value = 82700
$ws.Range("J90").Value = 82700
$ws.Range("L90").Value = 248100
$ws.Range("N90").Value = -260580
$ws.Range("H92").Value = 1067.6364
$ws.Range("I92").Value = 1141.6666
$ws.Range("K92").Value = 1141.6666
$ws.Range("M92").Value = 106.3334
$ws.Range("H112").Value = 4192.269
$ws.Range("J112").Value = 4319.96
$ws.Range("L112").Value = 12959.88
$ws.Range("N112").Value = -15175.88
$ws.Range("H116").Value = 26713636
$ws.Range("I116").Value = 20837208
$ws.Range("J116").Value = 38466492
$ws.Range("K116").Value = 20837208
$ws.Range("L116").Value = 38466492
$ws.Range("M116").Value = -20833766
$ws.Range("N116").Value = -38473376
$ws.Range("H125").Value = 3552
$ws.Range("I125").Value = 2292
$ws.Range("K125").Value = 20628
$ws.Range("M125").Value = -18168
$ws.Range("H132").Value = 148581.95
$ws.Range("I132").Value = 165920.44
$ws.Range("K132").Value = 497761.32
$ws.Range("M132").Value = -495231.32
$ws.Range("H135").Value = 41105.63
$ws.Range("I135").Value = 50088.145
$ws.Range("K135").Value = 450793.305
$ws.Range("M135").Value = -448258.305
$ws.Range("H137").Value = 4383.3
$ws.Range("I137").Value = 3439.6
$ws.Range("J137").Value = 4697.8667
$ws.Range("K137").Value = 10318.8
$ws.Range("L137").Value = 14093.6001
$ws.Range("M137").Value = -7768.799999999999
$ws.Range("N137").Value = -19193.6001
$ws.Range("H138").Value = 5712.9297
$ws.Range("I138").Value = 2239
$ws.Range("J138").Value = 5976.106
$ws.Range("K138").Value = 6717
$ws.Range("L138").Value = 17928.318
$ws.Range("M138").Value = -1577
$ws.Range("N138").Value = -28208.318

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3033.1177
$ws.Range("I32").Value = 1659.8276
$ws.Range("K32").Value = 1659.8276
$ws.Range("M32").Value = -1372.8276
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H36").Value = 5590.25
$ws.Range("I36").Value = 5590.25
$ws.Range("K36").Value = 5590.25
$ws.Range("M36").Value = -5244.25
$ws.Range("H45").Value = 2688.75
$ws.Range("I45").Value = 2616
$ws.Range("K45").Value = 2616
$ws.Range("M45").Value = -2239
$ws.Range("H61").Value = 11894.171
$ws.Range("I61").Value = 6530.387
$ws.Range("K61").Value = 6530.387
$ws.Range("M61").Value = -6318.387
$ws.Range("H74").Value = 4538.839
$ws.Range("I74").Value = 1813.6
$ws.Range("J74").Value = 7093.75
$ws.Range("K74").Value = 1813.6
$ws.Range("L74").Value = 7093.75
$ws.Range("M74").Value = -939.5999999999999
$ws.Range("N74").Value = -8841.75
$ws.Range("H77").Value = 4538.839
$ws.Range("I77").Value = 1813.6
$ws.Range("J77").Value = 7093.75
$ws.Range("K77").Value = 9068
$ws.Range("L77").Value = 35468.75
$ws.Range("M77").Value = -4700
$ws.Range("N77").Value = -44204.75
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H110").Value = 14743262
$ws.Range("I110").Value = 1277240.5
$ws.Range("K110").Value = 1277240.5
$ws.Range("M110").Value = -1275195.5
$ws.Range("H132").Value = 20806.559
$ws.Range("I132").Value = 21354.607
$ws.Range("J132").Value = 18249
$ws.Range("K132").Value = 64063.821
$ws.Range("L132").Value = 54747
$ws.Range("M132").Value = -61533.821
$ws.Range("N132").Value = -59807
$ws.Range("H136").Value = 11894.171
$ws.Range("I136").Value = 6530.387
$ws.Range("K136").Value = 19591.161
$ws.Range("M136").Value = -17041.161

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2731.6365
$ws.Range("J20").Value = 2738.3333
$ws.Range("L20").Value = 2738.3333
$ws.Range("N20").Value = -3232.3333
$ws.Range("H80").Value = 381.53845
$ws.Range("I80").Value = 690.5
$ws.Range("K80").Value = 690.5
$ws.Range("M80").Value = 307.5
$ws.Range("H83").Value = 381.53845
$ws.Range("I83").Value = 690.5
$ws.Range("K83").Value = 3452.5
$ws.Range("M83").Value = 1539.5
$ws.Range("H86").Value = 1993.3125
$ws.Range("I86").Value = 1662.625
$ws.Range("J86").Value = 2324
$ws.Range("K86").Value = 1662.625
$ws.Range("L86").Value = 2324
$ws.Range("M86").Value = -539.625
$ws.Range("N86").Value = -4570
$ws.Range("H89").Value = 1993.3125
$ws.Range("I89").Value = 1662.625
$ws.Range("J89").Value = 2324
$ws.Range("K89").Value = 8313.125
$ws.Range("L89").Value = 11620
$ws.Range("M89").Value = -2697.125
$ws.Range("N89").Value = -22852
$ws.Range("H94").Value = 548617.75
$ws.Range("I94").Value = 761685.8
$ws.Range("K94").Value = 761685.8
$ws.Range("M94").Value = -761234.8
$ws.Range("H107").Value = 2532.4666
$ws.Range("I107").Value = 2768.2307
$ws.Range("K107").Value = 2768.2307
$ws.Range("M107").Value = -848.2307000000001
$ws.Range("H109").Value = 84999.25
$ws.Range("J109").Value = 84999.25
$ws.Range("L109").Value = 84999.25
$ws.Range("N109").Value = -87773.25
$ws.Range("H134").Value = 3301.5557
$ws.Range("I134").Value = 1962.0526
$ws.Range("J134").Value = 6482.875
$ws.Range("K134").Value = 5886.1578
$ws.Range("L134").Value = 19448.625
$ws.Range("M134").Value = -3351.1578
$ws.Range("N134").Value = -24518.625

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 338.3125
$ws.Range("I7").Value = 367.9
$ws.Range("J7").Value = 289
$ws.Range("K7").Value = 367.9
$ws.Range("L7").Value = 289
$ws.Range("M7").Value = -254.9
$ws.Range("N7").Value = -515
$ws.Range("H22").Value = 325
$ws.Range("I22").Value = 333.33334
$ws.Range("K22").Value = 333.33334
$ws.Range("M22").Value = 16.66665999999998
$ws.Range("H31").Value = 22731922
$ws.Range("I31").Value = 45457000
$ws.Range("J31").Value = 6840.8184
$ws.Range("K31").Value = 45457000
$ws.Range("L31").Value = 6840.8184
$ws.Range("M31").Value = -45456705
$ws.Range("N31").Value = -7430.8184
$ws.Range("H34").Value = 22731922
$ws.Range("I34").Value = 45457000
$ws.Range("J34").Value = 6840.8184
$ws.Range("K34").Value = 45457000
$ws.Range("L34").Value = 6840.8184
$ws.Range("M34").Value = -45456798
$ws.Range("N34").Value = -7244.8184
$ws.Range("H55").Value = 75000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H58").Value = 1673254.4
$ws.Range("J58").Value = 9842.333000000001
$ws.Range("L58").Value = 9842.333000000001
$ws.Range("N58").Value = -10248.333
$ws.Range("H81").Value = 98600
$ws.Range("J81").Value = 98600
$ws.Range("L81").Value = 98600
$ws.Range("N81").Value = -100596
$ws.Range("H84").Value = 98600
$ws.Range("J84").Value = 98600
$ws.Range("L84").Value = 295800
$ws.Range("N84").Value = -305784
$ws.Range("H114").Value = 50991.668
$ws.Range("J114").Value = 50991.668
$ws.Range("L114").Value = 50991.668
$ws.Range("N114").Value = -59669.668
$ws.Range("H131").Value = 40795
$ws.Range("J131").Value = 40795
$ws.Range("L131").Value = 40795
$ws.Range("N131").Value = -50875
$ws.Range("H132").Value = 27216048
$ws.Range("I132").Value = 33336310
$ws.Range("K132").Value = 100008930
$ws.Range("M132").Value = -100006400
$ws.Range("H134").Value = 2284.889
$ws.Range("I134").Value = 2213.4119
$ws.Range("K134").Value = 6640.2357
$ws.Range("M134").Value = -4105.2357
$ws.Range("H136").Value = 1673254.4
$ws.Range("J136").Value = 9842.333000000001
$ws.Range("L136").Value = 29526.999
$ws.Range("N136").Value = -34626.999
$ws.Range("H140").Value = 70390
$ws.Range("J140").Value = 100780
$ws.Range("L140").Value = 100780
$ws.Range("N140").Value = -111140

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 44868.4
$ws.Range("I32").Value = 70615.664
$ws.Range("J32").Value = 6247.5
$ws.Range("K32").Value = 211846.992
$ws.Range("L32").Value = 18742.5
$ws.Range("M32").Value = -211563.992
$ws.Range("N32").Value = -19308.5
$ws.Range("H68").Value = 207468.84
$ws.Range("J68").Value = 327739.53
$ws.Range("L68").Value = 983218.5900000001
$ws.Range("N68").Value = -984840.5900000001
$ws.Range("H71").Value = 207468.84
$ws.Range("J71").Value = 327739.53
$ws.Range("L71").Value = 2949655.77
$ws.Range("N71").Value = -2957767.77
$ws.Range("H75").Value = 41671490
$ws.Range("I75").Value = 4750
$ws.Range("J75").Value = 50004836
$ws.Range("K75").Value = 14250
$ws.Range("L75").Value = 150014508
$ws.Range("M75").Value = -13252
$ws.Range("N75").Value = -150016504
$ws.Range("H78").Value = 41671490
$ws.Range("I78").Value = 4750
$ws.Range("J78").Value = 50004836
$ws.Range("K78").Value = 42750
$ws.Range("L78").Value = 450043524
$ws.Range("M78").Value = -37758
$ws.Range("N78").Value = -450053508
$ws.Range("H87").Value = 2500
$ws.Range("I87").Value = 2500
$ws.Range("K87").Value = 7500
$ws.Range("M87").Value = -6252
$ws.Range("H90").Value = 2500
$ws.Range("I90").Value = 2500
$ws.Range("K90").Value = 22500
$ws.Range("M90").Value = -16260
$ws.Range("H92").Value = 1040.5
$ws.Range("J92").Value = 1156.25
$ws.Range("L92").Value = 3468.75
$ws.Range("N92").Value = -5964.75
$ws.Range("H97").Value = 1115.2858
$ws.Range("I97").Value = 921.4
$ws.Range("K97").Value = 2764.2
$ws.Range("M97").Value = -2268.2
$ws.Range("H98").Value = 1627.3684
$ws.Range("I98").Value = 1050
$ws.Range("J98").Value = 1695.2941
$ws.Range("K98").Value = 3150
$ws.Range("L98").Value = 5085.8823
$ws.Range("M98").Value = -1652
$ws.Range("N98").Value = -8081.8823
$ws.Range("H113").Value = 540.12
$ws.Range("J113").Value = 657.7857
$ws.Range("L113").Value = 1973.3571
$ws.Range("N113").Value = -6313.3571
$ws.Range("H122").Value = 408.3
$ws.Range("I122").Value = 235.375
$ws.Range("K122").Value = 2118.375
$ws.Range("M122").Value = 331.625
$ws.Range("H129").Value = 2480.3
$ws.Range("I129").Value = 1325.8334
$ws.Range("J129").Value = 4212
$ws.Range("K129").Value = 3977.5002
$ws.Range("L129").Value = 12636
$ws.Range("M129").Value = 1022.4998
$ws.Range("N129").Value = -22636
$ws.Range("H131").Value = 15386722
$ws.Range("I131").Value = 66667212
$ws.Range("J131").Value = 6062997
$ws.Range("K131").Value = 200001636
$ws.Range("L131").Value = 18188991
$ws.Range("M131").Value = -199996596
$ws.Range("N131").Value = -18199071
$ws.Range("H134").Value = 4289.5264
$ws.Range("I134").Value = 3833.4
$ws.Range("K134").Value = 11500.2
$ws.Range("M134").Value = -6430.200000000001
$ws.Range("H137").Value = 5005016
$ws.Range("I137").Value = 1963.8889
$ws.Range("K137").Value = 5891.6667
$ws.Range("M137").Value = -791.6666999999998
$ws.Range("H140").Value = 5988
$ws.Range("I140").Value = 1552.7368
$ws.Range("J140").Value = 20033
$ws.Range("K140").Value = 4658.2104
$ws.Range("L140").Value = 60099
$ws.Range("M140").Value = 521.7896000000001
$ws.Range("N140").Value = -70459

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1139545.6
$ws.Range("I70").Value = 1705732.4
$ws.Range("J70").Value = 7172.0713
$ws.Range("K70").Value = 1705732.4
$ws.Range("L70").Value = 7172.0713
$ws.Range("M70").Value = -1705462.4
$ws.Range("N70").Value = -7712.0713
$ws.Range("H73").Value = 1139545.6
$ws.Range("I73").Value = 1705732.4
$ws.Range("J73").Value = 7172.0713
$ws.Range("K73").Value = 1705732.4
$ws.Range("L73").Value = 7172.0713
$ws.Range("M73").Value = -1704796.4
$ws.Range("N73").Value = -9044.0713
$ws.Range("H80").Value = 1522617.5
$ws.Range("I80").Value = 2084349.1
$ws.Range("J80").Value = 24666.334
$ws.Range("K80").Value = 2084349.1
$ws.Range("L80").Value = 24666.334
$ws.Range("M80").Value = -2083351.1
$ws.Range("N80").Value = -26662.334
$ws.Range("H83").Value = 1522617.5
$ws.Range("I83").Value = 2084349.1
$ws.Range("J83").Value = 24666.334
$ws.Range("K83").Value = 10421745.5
$ws.Range("L83").Value = 123331.67
$ws.Range("M83").Value = -10416753.5
$ws.Range("N83").Value = -133315.67
$ws.Range("H87").Value = 62000
$ws.Range("J87").Value = 62000
$ws.Range("L87").Value = 62000
$ws.Range("N87").Value = -64496
$ws.Range("H90").Value = 62000
$ws.Range("J90").Value = 62000
$ws.Range("L90").Value = 186000
$ws.Range("N90").Value = -198480
$ws.Range("H97").Value = 1846.9445
$ws.Range("I97").Value = 1458.88
$ws.Range("J97").Value = 2728.9092
$ws.Range("K97").Value = 1458.88
$ws.Range("L97").Value = 2728.9092
$ws.Range("M97").Value = -962.8800000000001
$ws.Range("N97").Value = -3720.9092
$ws.Range("H102").Value = 26319506
$ws.Range("I102").Value = 41668370
$ws.Range("K102").Value = 41668370
$ws.Range("M102").Value = -41666748
$ws.Range("H107").Value = 2507269.8
$ws.Range("I107").Value = 4763061
$ws.Range("K107").Value = 4763061
$ws.Range("M107").Value = -4761141
$ws.Range("H124").Value = 300000
$ws.Range("J124").Value = 300000
$ws.Range("L124").Value = 300000
$ws.Range("N124").Value = -309820
$ws.Range("H126").Value = 3510.1072
$ws.Range("I126").Value = 2094.7273
$ws.Range("K126").Value = 6284.1819
$ws.Range("M126").Value = -3814.1819
$ws.Range("H132").Value = 4502.372
$ws.Range("I132").Value = 4442.6484
$ws.Range("J132").Value = 4870.6665
$ws.Range("K132").Value = 13327.9452
$ws.Range("L132").Value = 14611.9995
$ws.Range("M132").Value = -10797.9452
$ws.Range("N132").Value = -19671.9995

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 92389400
$ws.Range("I16").Value = 120101520
$ws.Range("K16").Value = 120101520
$ws.Range("M16").Value = -120101350
$ws.Range("H22").Value = 1117.0667
$ws.Range("I22").Value = 1385.75
$ws.Range("K22").Value = 1385.75
$ws.Range("M22").Value = -1090.75
$ws.Range("H27").Value = 1117.0667
$ws.Range("I27").Value = 1385.75
$ws.Range("K27").Value = 1385.75
$ws.Range("M27").Value = -1278.75
$ws.Range("H68").Value = 5683310
$ws.Range("I68").Value = 5683310
$ws.Range("K68").Value = 5683310
$ws.Range("M68").Value = -5682561
$ws.Range("H71").Value = 5683310
$ws.Range("I71").Value = 5683310
$ws.Range("K71").Value = 28416550
$ws.Range("M71").Value = -28412806
$ws.Range("H93").Value = 1339.25
$ws.Range("I93").Value = 1267.2
$ws.Range("K93").Value = 1267.2
$ws.Range("M93").Value = -19.20000000000005
$ws.Range("H122").Value = 76927490
$ws.Range("I122").Value = 111114900
$ws.Range("J122").Value = 5826.25
$ws.Range("K122").Value = 333344700
$ws.Range("L122").Value = 17478.75
$ws.Range("M122").Value = -333342250
$ws.Range("N122").Value = -22378.75
$ws.Range("H132").Value = 4502.42
$ws.Range("I132").Value = 4299.61
$ws.Range("J132").Value = 5426.3335
$ws.Range("K132").Value = 12898.83
$ws.Range("L132").Value = 16279.0005
$ws.Range("M132").Value = -10368.83
$ws.Range("N132").Value = -21339.0005

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11267.143
$ws.Range("I74").Value = 8125
$ws.Range("J74").Value = 12524
$ws.Range("K74").Value = 8125
$ws.Range("L74").Value = 12524
$ws.Range("M74").Value = -7189
$ws.Range("N74").Value = -14396
$ws.Range("H77").Value = 11267.143
$ws.Range("I77").Value = 8125
$ws.Range("J77").Value = 12524
$ws.Range("K77").Value = 24375
$ws.Range("L77").Value = 37572
$ws.Range("M77").Value = -19695
$ws.Range("N77").Value = -46932
$ws.Range("H132").Value = 1327723.5
$ws.Range("I132").Value = 1639226.4
$ws.Range("K132").Value = 4917679.199999999
$ws.Range("M132").Value = -4915149.199999999
$ws.Range("H135").Value = 84333
$ws.Range("J135").Value = 84333
$ws.Range("L135").Value = 84333
$ws.Range("N135").Value = -94473
$ws.Range("H136").Value = 6667.921
$ws.Range("I136").Value = 4793.4443
$ws.Range("J136").Value = 8354.950000000001
$ws.Range("K136").Value = 14380.3329
$ws.Range("L136").Value = 25064.85
$ws.Range("M136").Value = -11830.3329
$ws.Range("N136").Value = -30164.85
